# Generate Report for Handoff
# Updates the "b.md" file's status/handoff info across the three sheets
# (Overview, zh-cn, de-de) to reflect that a new handoff was generated.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4a7190e2dfa2cb1f773c8186bfcfbee2dec56f4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b97a7bbd2b1ce9f151ac58a75fbe7998ab6eaab2/e2e/b.md."

# --- Overview sheet: row 3 corresponds to b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-05 04:42:38"

# --- zh-cn sheet: row 3 corresponds to b.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 04:42:34"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 corresponds to b.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 04:42:38"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
